# Add the Laubrock & Kliegl (2015) literature review entry as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "lit review" (the tab that's selected/active)

$row = 25

$citation = "Laubrock, J., & Kliegl, R. (2015). The eye-voice span during reading aloud. Frontiers in Psychology, 0. https://doi.org/10.3389/fpsyg.2015.01432"
$question = "What can the eye-voice span tell us about eye movement behavior during reading?  What factors affect the EVS?"
$summary = "Authors perform oral and silent reading tasks on the same set of German sentences (different participant groups) to perform exploratory analyses on the eye-voice span (EVS).  They find evidence that the oculomotor system is regulated by the cognitive system, with a relatively stable amount of information stored in (pre-articulatory) working memory. Given that this buffer is constantly updating during oral reading, online control is required, determining the when and where of `"corrective`" eye movements within the text."
$note1 = "Temporal and spatial EVS had lower interindividual variation than intraindividual variation; by offset ~254ms/9.7 letters, by onset ~561ms/16.2 letters."
$note2 = "Difficult words (low frequency or low predictability) require more processing and therefore lead to refixations, which attenuate the EVS."
$note3 = "When the EVS gets too large, two processes seek to correct the gap: refixations/longer fixations and, when fixation time is insufficient to control an expanding EVS, regressive eye movements."
$considerations = "By the time a participant begins processing the switch word, they will be articulating ~16 letters prior in the text.  This aligns well with the study design establishing the switch group as the switch word and two words prior/following and the preswitch group as the five words preceding the switch group, to ensure that behavioral effects related to the switch are seen either when the switch word is first fixated (at which time the voice will be articulating words in the preswitch group) or when articulation of the switch word is prepared/performed (switch group).  Given the assumptions about the working memory buffer here, it seems unlikely that behavioral effects will be seen very long after articulation of the switch word (that is, in the postswitch group)."
$quote = "`"...the overall pattern of results suggests that the EVS is quite flexible, and is adjusted according to cognitive, oculomotor, and articulatory demands. Given that the voice proceeds fairly linearly through the text, most of the adjustment is actually performed by the oculomotor system. The eyes, and also the mind, could in principle proceed faster than the voice, since silent reading is faster than oral reading. However, the eyes need to wait for the voice because the size of the working memory buffer is limited. The major target value in the system controlling the eyes during oral reading is a constant EVS at fixation offset of about 10 letters, translating into an average temporal EVS of about 560 ms...`""

# Write in the same order the shared-string table records them so new
# strings land at the expected indices (211..218).
$ws.Cells.Item($row, 1).Value = $citation        # A25 -> new string 211
$ws.Cells.Item($row, 4).Value = $note1            # D25 -> new string 212
$ws.Cells.Item($row, 5).Value = $note2            # E25 -> new string 213
$ws.Cells.Item($row, 3).Value = $summary          # C25 -> new string 214
$ws.Cells.Item($row, 6).Value = $note3            # F25 -> new string 215
$ws.Cells.Item($row, 7).Value = $considerations   # G25 -> new string 216
$ws.Cells.Item($row, 8).Value = $quote            # H25 -> new string 217
$ws.Cells.Item($row, 2).Value = $question         # B25 -> new string 218

# Match the row height Excel auto-computed for the other multi-line,
# wrapped-text rows (style 3 already carries wrapText + thin border).
$ws.Rows.Item($row).RowHeight = 176

$ws.Range("A3").Select()
